$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.102.23"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "1.848.52"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'0.6920"
$ws.Range("D6").Value = "'237.96"
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.07722"
$ws.Range("E8").Value = "  +8.81%  "
$ws.Range("D9").Value = "'0.3030"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "'23.22"
$ws.Range("E10").Value = "  -4.51%  "
$ws.Range("D11").Value = "'0.08109"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "1.855.01"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "'0.7234"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'5.205"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").Value = "'88.92"
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").Value = "29.116.52"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "'5.741"
$ws.Range("E17").Value = "  -4.40%  "
$ws.Range("D18").Value = "'0.000007779"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("D20").Value = "'235.43"
$ws.Range("E20").Value = "  -4.85%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.096.41"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'7.597"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").Value = "'8.962"
$ws.Range("D26").Value = "'160.74"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("E27").Value = "  -7.20%  "
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "'1.976"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'1.398"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("D31").Value = "'4.480"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").Value = "'4.014"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").Value = "'0.05217"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("D36").Value = "'1.024"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").Value = "'0.6986"
$ws.Range("E37").Value = "  -5.98%  "
$ws.Range("D38").Value = "'2.657"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("E39").Value = "  -4.13%  "
$ws.Range("D40").Value = "'2.675"
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("D41").Value = "'0.9173"
$ws.Range("E41").Value = "  +5.94%  "
$ws.Range("D42").Value = "'6.006"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "1.081.27"
$ws.Range("E43").Value = "  +3.86%  "
$ws.Range("D44").Value = "'0.4253"
$ws.Range("E44").Value = "  -4.26%  "
$ws.Range("D45").Value = "'70.41"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  +0.20%  "
$ws.Range("D47").Value = "'103.28"
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("D49").Value = "1.994.68"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Value = "'9.119"
$ws.Range("E50").Value = "  -4.00%  "
$ws.Range("D51").Value = "'6.978"
$ws.Range("E51").Value = "  -6.07%  "
